# Update crypto price/volume figures per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.855.73"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "2.018.94"
$ws.Range("E3").Value = "  -2.43%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'225.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.85%  "
$ws.Range("E6").Value = "  -3.71%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'54.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.86%  "
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.76%  "
$ws.Range("D12").Value = "2.318.46"
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("E13").Value = "  -4.28%  "
$ws.Range("D14").Value = "'20.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.77%  "
$ws.Range("E15").Value = "  -2.88%  "
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("D17").Value = "2.018.35"
$ws.Range("E17").Value = "  -2.76%  "
$ws.Range("D18").Value = "36.812.74"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").Value = "'6.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.53%  "
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("D21").Value = "0.0₃0818"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").Value = "'225.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  +3.15%  "
$ws.Range("D25").Value = "'2.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.92%  "
$ws.Range("D26").Value = "'165.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.80%  "
$ws.Range("D27").Value = "'9.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.55%  "
$ws.Range("E28").Value = "  -5.84%  "
$ws.Range("D29").Value = "'18.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.11%  "
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("E31").Value = "  -4.79%  "
$ws.Range("E32").Value = "  -3.49%  "
$ws.Range("D33").Value = "'0.0615"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.50%  "
$ws.Range("E34").Value = "  -4.49%  "
$ws.Range("E35").Value = "  -4.79%  "
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'3.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.96%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "1.488.15"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("E41").Value = "  -5.09%  "
$ws.Range("D42").Value = "'16.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("E43").Value = "  -2.66%  "
$ws.Range("D44").Value = "'94.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.91%  "
$ws.Range("E45").Value = "  -5.69%  "
$ws.Range("E46").Value = "  -5.52%  "
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("E48").Value = "  -3.79%  "
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").Value = "2.207.75"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").Value = "'3.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.01%  "
